$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1, J1 matching the style of the existing header row
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fill in the new numeric columns I and J for rows 2-8
$iValues = @(7, 2, 5, 8, 8, 5, 7)
$jValues = @(7, 2, 6, 9, 8, 5, 8)

for ($r = 2; $r -le 8; $r++) {
    $ws.Cells.Item($r, 9).Value = $iValues[$r - 2]
    $ws.Cells.Item($r, 10).Value = $jValues[$r - 2]
}
